$wb = $excel.ActiveWorkbook

# --- AVL sheet: "ships" (Passenger) lifetime value calibrated from 21 to 30 ---
$wsAVL = $wb.Worksheets.Item("AVL")
$wsAVL.Range("B6").Value = 30

# --- About sheet: add a "Notes:" row explaining the calibration ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A6").Value = "Notes:"
$wsAbout.Range("A6").Font.Bold = $true
$wsAbout.Range("B6").Value = "Passenger ships done through calibration, arriving at a value of 30"

# --- Selection on the About sheet moves down to B7 ---
[void]$wsAbout.Range("B7").Select()

# --- AVL becomes the active/selected tab, with its selection on B7 ---
$wsAVL.Activate()
[void]$wsAVL.Range("B7").Select()
